$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = [double]"0.9999999999999079"
$ws.Range("E2").Value = [double]"0.9999999999999079"

$ws.Range("D3").Value = [double]"1.989690385591352E-11"
$ws.Range("E3").Value = [double]"1.989690385591352E-11"

$ws.Range("D4").Value = [double]"1.257762884666898E-05"
$ws.Range("E4").Value = [double]"1.257762884666898E-05"

$ws.Range("D5").Value = [double]"3.482945881258352E-23"
$ws.Range("E5").Value = [double]"3.482945881258352E-23"

$ws.Range("D6").Value = [double]"7.301362456744138E-21"
$ws.Range("E6").Value = [double]"7.301362456744138E-21"

$ws.Range("D7").Value = [double]"0.9999958298017891"
$ws.Range("E7").Value = [double]"4.170198210862175E-06"

$ws.Range("D8").Value = [double]"0.9999998696958784"
$ws.Range("E8").Value = [double]"1.303041216038636E-07"

$ws.Range("D9").Value = [double]"0.9999999999985529"
$ws.Range("E9").Value = [double]"1.447064690296429E-12"

$ws.Range("D10").Value = [double]"0.9994718449871581"
$ws.Range("E10").Value = [double]"0.0005281550128418688"

$ws.Range("D11").Value = [double]"0.9999999999291127"
$ws.Range("E11").Value = [double]"7.08872960331064E-11"
$ws.Range("F11").Value = [double]"3.001678943634033"
